$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4220.690002768957,
    4220.690002768957,
    4167.734760263412,
    4167.734760263412,
    4165.596755217855,
    4165.596755217855,
    4165.596755217855,
    4165.596755217855,
    4165.596755217855,
    4143.277104190689,
    4143.277104190689
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
